$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$cellRef, [string]$val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextCell "D2" "69.662.58"
Set-TextCell "E2" "  +0.54%  "
Set-TextCell "D3" "3.709.25"
Set-TextCell "E3" "  +0.93%  "
Set-TextCell "D4" "0.999"
Set-TextCell "E4" "  -0.04%  "
Set-TextCell "D5" "673.17"
Set-TextCell "E5" "  -1.42%  "
Set-TextCell "D6" "162.01"
Set-TextCell "E6" "  +2.63%  "
Set-TextCell "E7" "  +0.00%  "
Set-TextCell "E8" "  +1.21%  "
Set-TextCell "E9" "  +0.93%  "
Set-TextCell "E10" "  +1.89%  "
Set-TextCell "E12" "  +1.45%  "
Set-TextCell "D13" "32.90"
Set-TextCell "E13" "  +2.36%  "
Set-TextCell "D14" "3.695.81"
Set-TextCell "E14" "  +0.19%  "
Set-TextCell "D15" "69.679.49"
Set-TextCell "E15" "  +0.53%  "
Set-TextCell "E16" "  +1.59%  "
Set-TextCell "D17" "16.30"
Set-TextCell "E17" "  +2.91%  "
Set-TextCell "D18" "6.52"
Set-TextCell "E18" "  +2.31%  "
Set-TextCell "D19" "473.44"
Set-TextCell "E19" "  +0.35%  "
Set-TextCell "D20" "9.82"
Set-TextCell "E20" "  -0.97%  "
Set-TextCell "E21" "  +1.10%  "
Set-TextCell "D22" "80.46"
Set-TextCell "E22" "  +0.67%  "
Set-TextCell "D23" "3.859.12"
Set-TextCell "E23" "  +1.03%  "
Set-TextCell "E24" "  +5.53%  "
Set-TextCell "E25" "  -0.03%  "
Set-TextCell "D26" "11.03"
Set-TextCell "E26" "  +1.15%  "
Set-TextCell "D27" "9.16"
Set-TextCell "E27" "  +0.47%  "
Set-TextCell "E28" "  -0.20%  "
Set-TextCell "E29" "  +0.50%  "
Set-TextCell "E30" "  +1.84%  "
Set-TextCell "B31" "NEARProtocol"
Set-TextCell "C31" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell "D31" "6.60"
Set-TextCell "E31" "  +0.91%  "
Set-TextCell "B32" "Kaspa"
Set-TextCell "C32" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell "D32" "0.168"
Set-TextCell "E32" "  +6.78%  "
Set-TextCell "E33" "  +0.02%  "
Set-TextCell "E34" "  +0.43%  "
Set-TextCell "D35" "3.698.58"
Set-TextCell "E35" "  +1.25%  "
Set-TextCell "D36" "8.56"
Set-TextCell "E36" "  +5.10%  "
Set-TextCell "D37" "6.12"
Set-TextCell "E37" "  +0.86%  "
Set-TextCell "E39" "  +2.46%  "
Set-TextCell "E40" "  -0.05%  "
Set-TextCell "E41" "  +2.06%  "
Set-TextCell "D42" "174.17"
Set-TextCell "E43" "  +0.22%  "
Set-TextCell "D44" "47.13"
Set-TextCell "E44" "  -0.86%  "
Set-TextCell "E45" "  +2.80%  "
Set-TextCell "E46" "  +1.85%  "
Set-TextCell "E47" "  +1.92%  "
Set-TextCell "D48" "27.76"
Set-TextCell "E48" "  +3.72%  "
Set-TextCell "E49" "  -0.16%  "
Set-TextCell "E50" "  +1.88%  "
Set-TextCell "E51" "  +1.61%  "
